$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ------------------------------------------------------------------
# 1) Shift the bottom block of the table down by one row (286-395 -> 287-396)
#    so a new "2024" year-header row can be inserted at row 286.
#    Copy bottom-up (395 -> 396 first) to avoid self-overlap corruption.
# ------------------------------------------------------------------
for ($r = 395; $r -ge 286; $r--) {
    $srcRange = $ws.Range("A" + $r + ":K" + $r)
    $dstRange = $ws.Range("A" + ($r + 1) + ":K" + ($r + 1))
    $srcRange.Copy($dstRange)
}

# ------------------------------------------------------------------
# 2) Build the new row 286 ("2024" year header) by copying the format
#    of the existing "2023" header row (273) - this reuses the exact
#    same style ids Excel already has (bold, centered, date-quote-prefixed).
# ------------------------------------------------------------------
$ws.Range("A273:K273").Copy($ws.Range("A286:K286"))
$ws.Range("A286").Value = "'2024"

# Row273's G column uses style 42, but the new 2024 header row instead
# reuses the earlier-style (13) used by older year headers - fix it up
# with a formats-only paste from a cell that already carries style 13.
$ws.Range("C283").Copy()
$ws.Range("G286").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3) Fill in the monthly SL/VL accrual + leave usage entries for
#    Jun-Nov 2023 (rows 279-284).
# ------------------------------------------------------------------
$ws.Range("C279").Value = 1.25
$ws.Range("C280").Value = 1.25
$ws.Range("C281").Value = 1.25
$ws.Range("C282").Value = 1.25

$ws.Range("K283").Value = "10/27-29/2023"
$ws.Range("B283").Value = "SL(3-0-0)"
$ws.Range("C283").Value = 1.25
$ws.Range("H283").Value = 3

$ws.Range("K284").Value = "12/9,10,16-18/2023"
$ws.Range("B284").Value = "VL(5-0-0)"
$ws.Range("C284").Value = 1.25
$ws.Range("D284").Value = 5

# ------------------------------------------------------------------
# 4) Grow the structured table by one row to cover the newly
#    inserted header row.
# ------------------------------------------------------------------
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A8:K396"))

# ------------------------------------------------------------------
# 5) Restore the view state (freeze-pane scroll position + selection)
#    as closely as the object model allows.
# ------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 269
$ws.Range("E285").Select()

$wb.Application.CalculateFull()
